# Fill column D (rows 1-8) on Sheet3 with "D1".."D8", mirroring the
# existing A/B/C columns, then leave the selection where the author
# left it (G16) on the already-active Sheet3 tab.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet3")
$ws.Activate()

for ($i = 1; $i -le 8; $i++) {
    $ws.Cells.Item($i, 4).Value = "D$i"
}

$ws.Range("G16").Select()
